$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.04834671834362325
$ws.Range("D2").Value = 0.102340726199543
$ws.Range("E2").Value = 0.08328565320068293
$ws.Range("F2").Value = 2.859129718111461
$ws.Range("G2").Value = 2.506027728869753
$ws.Range("H2").Value = 1.835817060205386
$ws.Range("J2").Value = 0.09912528355382477
$ws.Range("M2").Value = 1.969097753115591
$ws.Range("N2").Value = 1.573475684104238

$ws.Range("C3").Value = 0.0429180893500245
$ws.Range("D3").Value = 0.09562326158399514
$ws.Range("E3").Value = 0.08201667906204335
$ws.Range("F3").Value = 2.779587354382329
$ws.Range("G3").Value = 2.391787200681676
$ws.Range("H3").Value = 1.792288196531956
$ws.Range("J3").Value = 0.1004854461454627
$ws.Range("M3").Value = 1.797296184505669
$ws.Range("N3").Value = 1.46737723029517

$ws.Range("C4").Value = 0.03961074216603322
$ws.Range("D4").Value = 0.09146593579458795
$ws.Range("E4").Value = 0.08124143989931554
$ws.Range("F4").Value = 2.733217131338648
$ws.Range("G4").Value = 2.323906772236001
$ws.Range("H4").Value = 1.767051942773548
$ws.Range("J4").Value = 0.101373056645679
$ws.Range("M4").Value = 1.692114468396326
$ws.Range("N4").Value = 1.402568746280167

$ws.Range("C5").Value = 0.03826917861611889
$ws.Range("D5").Value = 0.08976330039014613
$ws.Range("E5").Value = 0.08092653530899341
$ws.Range("F5").Value = 2.714933397893304
$ws.Range("G5").Value = 2.29680413728849
$ws.Range("H5").Value = 1.757138278198653
$ws.Range("J5").Value = 0.1017479683468618
$ws.Range("M5").Value = 1.64932830999031
$ws.Range("N5").Value = 1.37624620687231

$ws.Range("C6").Value = 0.0380467794252155
$ws.Range("D6").Value = 0.08948005900109024
$ws.Range("E6").Value = 0.08087430753906766
$ws.Range("F6").Value = 2.711934150037223
$ws.Range("G6").Value = 2.292337248872542
$ws.Range("H6").Value = 1.75551435721249
$ws.Range("J6").Value = 0.1018110199009605
$ws.Range("M6").Value = 1.642228293926479
$ws.Range("N6").Value = 1.371880735857076

$ws.Range("C7").Value = 0.03959262462196023
$ws.Range("D7").Value = 0.09144300816663531
$ws.Range("E7").Value = 0.08123718885371289
$ws.Range("F7").Value = 2.732968081431096
$ws.Range("G7").Value = 2.323539005616965
$ws.Range("H7").Value = 1.766916750145924
$ws.Range("J7").Value = 0.1013780593621387
$ws.Range("M7").Value = 1.691537131158313
$ws.Range("N7").Value = 1.402213393417071

$ws.Range("C8").Value = 0.04646938155842406
$ws.Range("D8").Value = 0.1000311936732601
$ws.Range("E8").Value = 0.08284731382360633
$ws.Range("F8").Value = 2.831185936632068
$ws.Range("G8").Value = 2.466161629007786
$ws.Range("H8").Value = 1.820496161478729
$ws.Range("J8").Value = 0.09958338478604745
$ws.Range("M8").Value = 1.909796864458301
$ws.Range("N8").Value = 1.536824775533688

$ws.Range("C9").Value = 0.06017388422064585
$ws.Range("D9").Value = 0.1166246457661231
$ws.Range("E9").Value = 0.08603482341689528
$ws.Range("F9").Value = 3.043778462150698
$ws.Range("G9").Value = 2.764281490661233
$ws.Range("H9").Value = 1.937607631688991
$ws.Range("J9").Value = 0.09647983175034369
$ws.Range("M9").Value = 2.340273983668368
$ws.Range("N9").Value = 1.803357838202118

$ws.Range("C10").Value = 0.07039690942414722
$ws.Range("D10").Value = 0.1286826853153542
$ws.Range("E10").Value = 0.08839390312727602
$ws.Range("F10").Value = 3.212731734877281
$ws.Range("G10").Value = 2.995242513984238
$ws.Range("H10").Value = 2.03130429728634
$ws.Range("J10").Value = 0.09445242960935119
$ws.Range("M10").Value = 2.658159702887247
$ws.Range("N10").Value = 2.000619432434746

$ws.Range("C11").Value = 0.0750857760916972
$ws.Range("D11").Value = 0.1341433251698163
$ws.Range("E11").Value = 0.08947063217429374
$ws.Range("F11").Value = 3.292493132229595
$ws.Range("G11").Value = 3.103058625009965
$ws.Range("H11").Value = 2.075661873806951
$ws.Range("J11").Value = 0.09358485403232919
$ws.Range("M11").Value = 2.803152260169782
$ws.Range("N11").Value = 2.090648376778176

$ws.Range("C12").Value = 0.07686721218748005
$ws.Range("D12").Value = 0.1362079028945402
$ws.Range("E12").Value = 0.08987885200702195
$ws.Range("F12").Value = 3.323124704485338
$ws.Range("G12").Value = 3.144293845121695
$ws.Range("H12").Value = 2.092714045495597
$ws.Range("J12").Value = 0.09326418212691223
$ws.Range("M12").Value = 2.858114190227809
$ws.Range("N12").Value = 2.12477993992286

$ws.Range("C13").Value = 0.0764832816204688
$ws.Range("D13").Value = 0.1357633988632045
$ws.Range("E13").Value = 0.08979091338918366
$ws.Range("F13").Value = 3.31650846326491
$ws.Range("G13").Value = 3.135394779832836
$ws.Range("H13").Value = 2.089030128451839
$ws.Range("J13").Value = 0.09333289506385256
$ws.Range("M13").Value = 2.846274628535383
$ws.Range("N13").Value = 2.117427379628168

$ws.Range("C14").Value = 0.07523221636489552
$ws.Range("D14").Value = 0.1343132426142546
$ws.Range("E14").Value = 0.0895042070847154
$ws.Range("F14").Value = 3.295004580331181
$ws.Range("G14").Value = 3.106442832510822
$ws.Range("H14").Value = 2.077059625249717
$ws.Range("J14").Value = 0.093558314664552
$ws.Range("M14").Value = 2.807672874384252
$ws.Range("N14").Value = 2.093455623312877

$ws.Range("C15").Value = 0.07446667690845743
$ws.Range("D15").Value = 0.1334245658594853
$ws.Range("E15").Value = 0.0893286536581499
$ws.Range("F15").Value = 3.281888829673392
$ws.Range("G15").Value = 3.08876237242788
$ws.Range("H15").Value = 2.069760710178628
$ws.Range("J15").Value = 0.0936974141374094
$ws.Range("M15").Value = 2.784035563728111
$ws.Range("N15").Value = 2.078777303659081

$ws.Range("C16").Value = 0.07009128242640372
$ws.Range("D16").Value = 0.128325349032707
$ws.Range("E16").Value = 0.08832360584973031
$ws.Range("F16").Value = 3.207578428866498
$ws.Range("G16").Value = 2.988252853026495
$ws.Range("H16").Value = 2.028440794965434
$ws.Range("J16").Value = 0.09451022792665853
$ws.Range("M16").Value = 2.648691928735474
$ws.Range("N16").Value = 1.994741508298034

$ws.Range("C17").Value = 0.06741718970188515
$ws.Range("D17").Value = 0.1251910440753221
$ws.Range("E17").Value = 0.08770793658483456
$ws.Range("F17").Value = 3.162742232194802
$ws.Range("G17").Value = 2.927306536457195
$ws.Range("H17").Value = 2.003540482019332
$ws.Range("J17").Value = 0.09502286870354126
$ws.Range("M17").Value = 2.565762254319139
$ws.Range("N17").Value = 1.943261558554667

$ws.Range("C18").Value = 0.0658827117388654
$ws.Range("D18").Value = 0.1233859391956713
$ws.Range("E18").Value = 0.08735415774386723
$ws.Range("F18").Value = 3.137226328917507
$ws.Range("G18").Value = 2.892510006736359
$ws.Range("H18").Value = 1.98938140202182
$ws.Range("J18").Value = 0.09532287504822889
$ws.Range("M18").Value = 2.518099441273478
$ws.Range("N18").Value = 1.913679448494747

$ws.Range("C19").Value = 0.06536376989882342
$ws.Range("D19").Value = 0.1227743517713549
$ws.Range("E19").Value = 0.08723443318119095
$ws.Range("F19").Value = 3.128633613257136
$ws.Range("G19").Value = 2.880772480045778
$ws.Range("H19").Value = 1.984615197284143
$ws.Range("J19").Value = 0.09542533654951768
$ws.Range("M19").Value = 2.501967808823764
$ws.Range("N19").Value = 1.903668310778158

$ws.Range("C20").Value = 0.06770147745058352
$ws.Range("D20").Value = 0.1255249358148802
$ws.Range("E20").Value = 0.08777344086608352
$ws.Range("F20").Value = 3.167486826902319
$ws.Range("G20").Value = 2.933767568992835
$ws.Range("H20").Value = 2.006174261683157
$ws.Range("J20").Value = 0.09496776436490961
$ws.Range("M20").Value = 2.574586511826055
$ws.Range("N20").Value = 1.948738829943295

$ws.Range("C21").Value = 0.07559952259620673
$ws.Range("D21").Value = 0.1347392744398803
$ws.Range("E21").Value = 0.08958840670807078
$ws.Range("F21").Value = 3.301309103563426
$ws.Range("G21").Value = 3.114935561753498
$ws.Range("H21").Value = 2.080568689203801
$ws.Range("J21").Value = 0.09349189018511694
$ws.Range("M21").Value = 2.819009605557312
$ws.Range("N21").Value = 2.100495656630414

$ws.Range("C22").Value = 0.0807956992394594
$ws.Range("D22").Value = 0.1407425863659313
$ws.Range("E22").Value = 0.09077741145046048
$ws.Range("F22").Value = 3.391267309600551
$ws.Range("G22").Value = 3.235720910024611
$ws.Range("H22").Value = 2.130678271325053
$ws.Range("J22").Value = 0.09257313355391972
$ws.Range("M22").Value = 2.979083535136965
$ws.Range("N22").Value = 2.199906887718782

$ws.Range("C23").Value = 0.07801914208826588
$ws.Range("D23").Value = 0.13754012652754
$ws.Range("E23").Value = 0.09014256831390455
$ws.Range("F23").Value = 3.343022984556825
$ws.Range("G23").Value = 3.171033527339887
$ws.Range("H23").Value = 2.103795789901142
$ws.Range("J23").Value = 0.09305930120613226
$ws.Range("M23").Value = 2.893618546225554
$ws.Range("N23").Value = 2.146829119470851

$ws.Range("C24").Value = 0.06757294197278441
$ws.Range("D24").Value = 0.1253739931958364
$ws.Range("E24").Value = 0.08774382582514306
$ws.Range("F24").Value = 3.1653409833674
$ws.Range("G24").Value = 2.93084578231975
$ws.Range("H24").Value = 2.004983043031132
$ws.Range("J24").Value = 0.09499266058269384
$ws.Range("M24").Value = 2.570597019729718
$ws.Range("N24").Value = 1.946262510495757

$ws.Range("C25").Value = 0.0564408731712831
$ws.Range("D25").Value = 0.112160386338914
$ws.Range("E25").Value = 0.08516941510098164
$ws.Range("F25").Value = 2.984067929782498
$ws.Range("G25").Value = 2.681587665541372
$ws.Range("H25").Value = 1.904604694182638
$ws.Range("J25").Value = 0.09727498471777452
$ws.Range("M25").Value = 2.223543303372224
$ws.Range("N25").Value = 1.730992261561681
